# Apply the edit described by the diff:
#  - Add a new "N" column (column B) to the results table, shifting the
#    existing Promedio/Desv. Est./Mín/Máx./Diferencia de medias columns
#    one slot to the right (C..G) for the header (row 3) and data rows
#    (5..7). Row 4, the footnote-marker row, keeps its existing (1)-(6)
#    values in place and simply gains a new (7) marker in column G.
#  - Populate the new N column with the sample sizes (as real numbers).
#  - Refresh the KNN / XGBoost statistics with the new modelling results.
#  - Update the explanatory note under the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, [string]$text, [string]$fmtFrom)
    # Force Excel to store the value as text even when it looks like a
    # number (e.g. "0.0382"), then copy the direct formatting from a
    # donor cell (left untouched by this script) so the style index
    # matches exactly (font/size/etc.) instead of drifting to a new one.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range($fmtFrom).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# --- Column structure -------------------------------------------------
# New column G takes over the width previously used by the other data
# columns (COM value 13.17 renders to the stored width "14").
$ws.Columns.Item(7).ColumnWidth = 13.17

# --- Row 3: header (shift Promedio..Diferencia de medias right by one) --
# Donor = A3 (untouched, already styled like the rest of the header row).
Set-TextValue "B3" "N" "A3"
Set-TextValue "C3" "Promedio" "A3"
Set-TextValue "D3" "Desv. Est." "A3"
Set-TextValue "E3" "Mín" "A3"
Set-TextValue "F3" "Máx." "A3"
Set-TextValue "G3" "Diferencia de medias" "A3"

# --- Row 4: footnote markers ---------------------------------------------
# (1)-(6) already sit in A4:F4 and stay exactly where they are; only the
# new (7) marker for the added column is appended. Donor = A4.
Set-TextValue "G4" "(7)" "A4"

# --- Row 5: Datos originales ----------------------------------------------
# Donor = A5.
Set-TextValue "C5" "0.0001" "A5"
Set-TextValue "D5" "1.0164" "A5"
Set-TextValue "E5" "-5.1993" "A5"
Set-TextValue "F5" "5.1993" "A5"
Set-TextValue "G5" "—" "A5"
$ws.Range("B5").Value = 1236

# --- Row 6: KNN -------------------------------------------------------------
# Donor = A6.
Set-TextValue "C6" "0.0382" "A6"
Set-TextValue "D6" "0.6675" "A6"
Set-TextValue "E6" "-1.7669" "A6"
Set-TextValue "F6" "2.5278" "A6"
Set-TextValue "G6" "0.0381" "A6"
$ws.Range("B6").Value = 5000

# --- Row 7: XGBoost ----------------------------------------------------------
# Donor = A7.
Set-TextValue "C7" "0.0851" "A7"
Set-TextValue "D7" "0.7715" "A7"
Set-TextValue "E7" "-2.3997" "A7"
Set-TextValue "F7" "3.4618" "A7"
Set-TextValue "G7" "0.0851" "A7"
$ws.Range("B7").Value = 5000

# --- Merged title / note ranges now span A:G --------------------------------
$ws.Range("A1:F1").UnMerge()
$ws.Range("A1:G1").Merge()
$ws.Range("A10:F10").UnMerge()
$ws.Range("A10:G10").Merge()

# --- Explanatory note ---------------------------------------------------------
$ws.Range("A10").Value = "(a) Elaboración propia en base a datos procesados (cluster con nscore y predicciones de modelos KNN y XGBoost). (b) N = número de datos. Desv. Est. = desviación estándar. Mín = mínimo valor observado. Máx. = máximo valor observado. Diferencia de medias = diferencia respecto a la media de datos originales."
